$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update the Runmode column (C2:C12) from "Y" to "N" so only the D suite runs
$ws.Range("C2:C12").Value = "N"

# Reflect the new selection (activeCell=C2, sqref=C2:C12) used while making this change
$ws.Range("C2:C12").Select()
